# Weekly update: insert a new daily record at row 43 (pushing existing
# history down by one row) and populate it with the newest data point.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 43; Excel shifts rows 43:162 down to 44:163
# and copies formatting (including the date style on column D) from the
# row that is being pushed down.
$ws.Rows("43:43").Insert()

# Populate the newly inserted row with the new observation.
$ws.Cells.Item(43, 1).Value2 = 8
$ws.Cells.Item(43, 2).Value2 = "Terminal La Palmera de La Serena"
$ws.Cells.Item(43, 3).Value2 = "Coquimbo"
$ws.Cells.Item(43, 4).Value2 = 44414
$ws.Cells.Item(43, 5).Value2 = 4
$ws.Cells.Item(43, 6).Value2 = 100114013
$ws.Cells.Item(43, 7).Value2 = "Zanahoria"
$ws.Cells.Item(43, 8).Value2 = "Sin especificar"
$ws.Cells.Item(43, 9).Value2 = "Primera"
$ws.Cells.Item(43, 10).Value2 = 700
$ws.Cells.Item(43, 11).Value2 = 5000
$ws.Cells.Item(43, 12).Value2 = 5500
$ws.Cells.Item(43, 13).Value2 = 5250
$ws.Cells.Item(43, 14).Value2 = "$/saco 20 kilos"
$ws.Cells.Item(43, 15).Value2 = "Provincia del Elquí"
$ws.Cells.Item(43, 16).Value2 = 262
$ws.Cells.Item(43, 17).Value2 = 20
$ws.Cells.Item(43, 18).Value2 = "Hortaliza"
